$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1276.711
$ws.Range("I70").Value = 930
$ws.Range("J70").Value = 1375.7715
$ws.Range("K70").Value = 2790
$ws.Range("L70").Value = 4127.3145
$ws.Range("M70").Value = -2520
$ws.Range("N70").Value = -4667.3145
$ws.Range("H73").Value = 1276.711
$ws.Range("I73").Value = 930
$ws.Range("J73").Value = 1375.7715
$ws.Range("K73").Value = 2790
$ws.Range("L73").Value = 4127.3145
$ws.Range("M73").Value = -1854
$ws.Range("N73").Value = -5999.3145
$ws.Range("H80").Value = 478.9
$ws.Range("I80").Value = 200
$ws.Range("J80").Value = 548.625
$ws.Range("K80").Value = 600
$ws.Range("L80").Value = 1645.875
$ws.Range("M80").Value = 398
$ws.Range("N80").Value = -3641.875
$ws.Range("H83").Value = 478.9
$ws.Range("I83").Value = 200
$ws.Range("J83").Value = 548.625
$ws.Range("K83").Value = 1800
$ws.Range("L83").Value = 4937.625
$ws.Range("M83").Value = 3192
$ws.Range("N83").Value = -14921.625
$ws.Range("H100").Value = 22224624
$ws.Range("I100").Value = 34190360
$ws.Range("J100").Value = 2544.5715
$ws.Range("K100").Value = 34190360
$ws.Range("L100").Value = 2544.5715
$ws.Range("M100").Value = -34189819
$ws.Range("N100").Value = -3626.5715
$ws.Range("H113").Value = 858585.4399999999
$ws.Range("J113").Value = 4350
$ws.Range("L113").Value = 4350
$ws.Range("N113").Value = -10858
$ws.Range("H132").Value = 4322497
$ws.Range("I132").Value = 7879248.5
$ws.Range("J132").Value = 3584.5
$ws.Range("K132").Value = 23637745.5
$ws.Range("L132").Value = 10753.5
$ws.Range("M132").Value = -23635215.5
$ws.Range("N132").Value = -15813.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10884
$ws.Range("I32").Value = 7717.241
$ws.Range("J32").Value = 32076.924
$ws.Range("K32").Value = 7717.241
$ws.Range("L32").Value = 32076.924
$ws.Range("M32").Value = -7430.241
$ws.Range("N32").Value = -32650.924

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 46433.332
$ws.Range("J18").Value = 46433.332
$ws.Range("L18").Value = 46433.332
$ws.Range("N18").Value = -46893.332
$ws.Range("H99").Value = 7837.2
$ws.Range("I99").Value = 10523.846
$ws.Range("J99").Value = 2847.7144
$ws.Range("K99").Value = 10523.846
$ws.Range("L99").Value = 2847.7144
$ws.Range("M99").Value = -9025.846
$ws.Range("N99").Value = -5843.7144
$ws.Range("H117").Value = 42178
$ws.Range("J117").Value = 42178
$ws.Range("L117").Value = 42178
$ws.Range("N117").Value = -51356
$ws.Range("H126").Value = 7837.2
$ws.Range("I126").Value = 10523.846
$ws.Range("J126").Value = 2847.7144
$ws.Range("K126").Value = 31571.538
$ws.Range("L126").Value = 8543.143199999999
$ws.Range("M126").Value = -29101.538
$ws.Range("N126").Value = -13483.1432

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 250.26666
$ws.Range("I2").Value = 60
$ws.Range("K2").Value = 360
$ws.Range("M2").Value = -247
$ws.Range("H22").Value = 1600
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1966.6666
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 5899.9998
$ws.Range("M22").Value = -1331
$ws.Range("N22").Value = -6237.9998
$ws.Range("H27").Value = 1600
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1966.6666
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 5899.9998
$ws.Range("M27").Value = -1398
$ws.Range("N27").Value = -6103.9998
$ws.Range("H106").Value = 6197.143
$ws.Range("J106").Value = 6197.143
$ws.Range("L106").Value = 18591.429
$ws.Range("N106").Value = -20483.429
$ws.Range("H107").Value = 778.6875
$ws.Range("I107").Value = 632.5
$ws.Range("J107").Value = 827.4167
$ws.Range("K107").Value = 1897.5
$ws.Range("L107").Value = 2482.2501
$ws.Range("M107").Value = 22.5
$ws.Range("N107").Value = -6322.2501
$ws.Range("H113").Value = 722013.75
$ws.Range("J113").Value = 536.875
$ws.Range("L113").Value = 1610.625
$ws.Range("N113").Value = -5950.625
$ws.Range("H120").Value = 14998.214
$ws.Range("I120").Value = 8999.333000000001
$ws.Range("J120").Value = 19497.375
$ws.Range("K120").Value = 26997.999
$ws.Range("L120").Value = 58492.125
$ws.Range("M120").Value = -22159.999
$ws.Range("N120").Value = -68168.125
$ws.Range("H122").Value = 864.8
$ws.Range("I122").Value = 708.6875
$ws.Range("J122").Value = 1043.2142
$ws.Range("K122").Value = 6378.1875
$ws.Range("L122").Value = 9388.927799999999
$ws.Range("M122").Value = -3928.1875
$ws.Range("N122").Value = -14288.9278
$ws.Range("H131").Value = 49830.145
$ws.Range("I131").Value = 273.75
$ws.Range("J131").Value = 61843.816
$ws.Range("K131").Value = 821.25
$ws.Range("L131").Value = 185531.448
$ws.Range("M131").Value = 4218.75
$ws.Range("N131").Value = -195611.448

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 48850
$ws.Range("J108").Value = 48850
$ws.Range("L108").Value = 48850
$ws.Range("N108").Value = -56530
$ws.Range("H126").Value = 2312.5
$ws.Range("I126").Value = 2150
$ws.Range("J126").Value = 2475
$ws.Range("K126").Value = 6450
$ws.Range("L126").Value = 7425
$ws.Range("M126").Value = -3980
$ws.Range("N126").Value = -12365

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2276.875
$ws.Range("I61").Value = 1406.6666
$ws.Range("J61").Value = 3147.0833
$ws.Range("K61").Value = 1406.6666
$ws.Range("L61").Value = 3147.0833
$ws.Range("M61").Value = -1204.6666
$ws.Range("N61").Value = -3551.0833
$ws.Range("H82").Value = 3262.5
$ws.Range("I82").Value = 3033.3333
$ws.Range("J82").Value = 3400
$ws.Range("K82").Value = 3033.3333
$ws.Range("L82").Value = 3400
$ws.Range("M82").Value = -2672.3333
$ws.Range("N82").Value = -4122
$ws.Range("H85").Value = 3262.5
$ws.Range("I85").Value = 3033.3333
$ws.Range("J85").Value = 3400
$ws.Range("K85").Value = 3033.3333
$ws.Range("L85").Value = 3400
$ws.Range("M85").Value = -1785.3333
$ws.Range("N85").Value = -5896
$ws.Range("H113").Value = 2276.875
$ws.Range("I113").Value = 1406.6666
$ws.Range("J113").Value = 3147.0833
$ws.Range("K113").Value = 1406.6666
$ws.Range("L113").Value = 3147.0833
$ws.Range("M113").Value = 763.3334
$ws.Range("N113").Value = -7487.0833

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 15152531
$ws.Range("I100").Value = 90909090
$ws.Range("J100").Value = 1220
$ws.Range("K100").Value = 181818180
$ws.Range("L100").Value = 2440
$ws.Range("M100").Value = -181817639
$ws.Range("N100").Value = -3522
